# Runtime update for Injuries_Master_Clubs workbook (2025-10-29 21:09:22)
# - "Фарранс Дэвид" (СИБ / Сибирь) is no longer injured -> removed from the
#   "snapshot" sheet and appended to the "returned" sheet.
# - The "new_injured" sheet's pending entry (Профака Лука) has already been
#   folded into the snapshot, so it is cleared back down to just the header.
# - All remaining "snapshot" rows get a refreshed scraped_at timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) snapshot: remove the row for Фарранс Дэвид (СИБ), shifting rows below up
# ---------------------------------------------------------------------------
$snapshot = $wb.Worksheets.Item("snapshot")

$returnedTeamAbbr   = "СИБ"
$returnedTeamName   = "Сибирь"
$returnedPlayerName = "Фарранс Дэвид"
$returnedPlayerUid  = "1369_СИБ_фаррансдэвид"

$targetRow = 0
for ($r = 2; $r -le $snapshot.UsedRange.Rows.Count + 1; $r++) {
    if ($snapshot.Cells.Item($r, 1).Value2 -eq $returnedTeamAbbr -and $snapshot.Cells.Item($r, 4).Value2 -eq $returnedPlayerName) {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    $snapshot.Rows.Item($targetRow).Delete()
}

# ---------------------------------------------------------------------------
# 2) snapshot: refresh the scraped_at timestamp (column K) for every
#    remaining data row
# ---------------------------------------------------------------------------
$scrapedAt = @(
    "2025-10-29T13:08:08.426764+00:00",
    "2025-10-29T13:08:10.467424+00:00",
    "2025-10-29T13:08:10.467441+00:00",
    "2025-10-29T13:08:10.467450+00:00",
    "2025-10-29T13:08:10.467458+00:00",
    "2025-10-29T13:08:10.467465+00:00",
    "2025-10-29T13:08:12.568705+00:00",
    "2025-10-29T13:08:12.568721+00:00",
    "2025-10-29T13:08:14.598217+00:00",
    "2025-10-29T13:08:16.577779+00:00",
    "2025-10-29T13:08:16.577795+00:00",
    "2025-10-29T13:08:18.800987+00:00",
    "2025-10-29T13:08:18.801002+00:00",
    "2025-10-29T13:08:18.801010+00:00",
    "2025-10-29T13:08:18.801017+00:00",
    "2025-10-29T13:08:24.982635+00:00",
    "2025-10-29T13:08:26.960682+00:00",
    "2025-10-29T13:08:29.006311+00:00",
    "2025-10-29T13:08:30.989093+00:00",
    "2025-10-29T13:08:30.989110+00:00",
    "2025-10-29T13:08:30.989118+00:00",
    "2025-10-29T13:08:33.448093+00:00",
    "2025-10-29T13:08:33.448109+00:00",
    "2025-10-29T13:08:33.448116+00:00",
    "2025-10-29T13:08:33.448123+00:00",
    "2025-10-29T13:08:33.448131+00:00",
    "2025-10-29T13:08:40.431219+00:00",
    "2025-10-29T13:08:40.431235+00:00",
    "2025-10-29T13:08:40.431243+00:00",
    "2025-10-29T13:08:40.431251+00:00",
    "2025-10-29T13:08:42.383583+00:00",
    "2025-10-29T13:08:42.383600+00:00",
    "2025-10-29T13:08:42.383608+00:00",
    "2025-10-29T13:08:44.394858+00:00",
    "2025-10-29T13:08:44.394877+00:00",
    "2025-10-29T13:08:44.394885+00:00",
    "2025-10-29T13:08:44.394892+00:00",
    "2025-10-29T13:08:44.394903+00:00",
    "2025-10-29T13:08:44.394910+00:00",
    "2025-10-29T13:08:44.394919+00:00",
    "2025-10-29T13:08:44.394927+00:00",
    "2025-10-29T13:08:44.394934+00:00",
    "2025-10-29T13:08:46.383149+00:00",
    "2025-10-29T13:08:46.383167+00:00",
    "2025-10-29T13:08:50.456984+00:00",
    "2025-10-29T13:08:52.858587+00:00",
    "2025-10-29T13:08:52.858603+00:00",
    "2025-10-29T13:08:52.858610+00:00",
    "2025-10-29T13:08:52.858617+00:00"
)

for ($i = 0; $i -lt $scrapedAt.Length; $i++) {
    $row = $i + 2
    $snapshot.Cells.Item($row, 11).Value = $scrapedAt[$i]
}

# ---------------------------------------------------------------------------
# 3) returned: append the Фарранс Дэвид "RETURN" record
# ---------------------------------------------------------------------------
$returned = $wb.Worksheets.Item("returned")
$newRow = $returned.UsedRange.Rows.Count + 1

$returned.Cells.Item($newRow, 1).Value = $returnedTeamAbbr
$returned.Cells.Item($newRow, 2).Value = $returnedTeamName
$returned.Cells.Item($newRow, 3).Value = $returnedPlayerName
$returned.Cells.Item($newRow, 4).Value = $returnedPlayerUid
$returned.Cells.Item($newRow, 5).Value = "RETURN"
$returned.Cells.Item($newRow, 6).Value = "2025-10-29T21:08:53.360093+08:00"
$returned.Cells.Item($newRow, 7).Value = "'2025-10-29"

# ---------------------------------------------------------------------------
# 4) new_injured: the pending Профака Лука entry already exists in the
#    snapshot, so clear the sheet back down to just the header row
# ---------------------------------------------------------------------------
$newInjured = $wb.Worksheets.Item("new_injured")
if ($newInjured.UsedRange.Rows.Count -ge 2) {
    $newInjured.Rows.Item(2).Delete()
}

Write-Host "edit complete"
